{"js": "// Replace each three-digit-divided-by-one-digit expression in the\n// document's table cells with its new value, per the commit diff.\n// Old -> New text pairs (each old value is unique in the document).\nconst replacements = [\n  [\"528\u00f75=\", \"379\u00f76=\"],\n  [\"421\u00f77=\", \"888\u00f72=\"],\n  [\"840\u00f75=\", \"665\u00f75=\"],\n  [\"255\u00f73=\", \"474\u00f74=\"],\n  [\"541\u00f73=\", \"117\u00f79=\"],\n  [\"229\u00f79=\", \"805\u00f77=\"],\n  [\"688\u00f74=\", \"854\u00f76=\"],\n  [\"242\u00f73=\", \"235\u00f79=\"],\n  [\"700\u00f75=\", \"511\u00f73=\"],\n  [\"778\u00f78=\", \"267\u00f75=\"],\n  [\"702\u00f76=\", \"525\u00f74=\"],\n  [\"102\u00f78=\", \"347\u00f79=\"],\n  [\"289\u00f72=\", \"245\u00f72=\"],\n  [\"404\u00f75=\", \"422\u00f76=\"],\n  [\"947\u00f78=\", \"314\u00f74=\"],\n  [\"214\u00f78=\", \"556\u00f76=\"],\n  [\"491\u00f73=\", \"324\u00f72=\"],\n  [\"685\u00f75=\", \"222\u00f79=\"],\n  [\"960\u00f76=\", \"431\u00f78=\"],\n  [\"715\u00f76=\", \"462\u00f77=\"],\n  [\"868\u00f72=\", \"401\u00f76=\"],\n  [\"883\u00f74=\", \"355\u00f77=\"],\n  [\"438\u00f75=\", \"839\u00f76=\"],\n  [\"223\u00f76=\", \"503\u00f74=\"],\n  [\"713\u00f74=\", \"659\u00f72=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit-divided-by-one-digit expression in the\n# document's table cells with its new value, per the commit diff.\n# Old -> New text pairs (each old value is unique in the document).\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @{ Old = \"528\u00f75=\"; New = \"379\u00f76=\" },\n  @{ Old = \"421\u00f77=\"; New = \"888\u00f72=\" },\n  @{ Old = \"840\u00f75=\"; New = \"665\u00f75=\" },\n  @{ Old = \"255\u00f73=\"; New = \"474\u00f74=\" },\n  @{ Old = \"541\u00f73=\"; New = \"117\u00f79=\" },\n  @{ Old = \"229\u00f79=\"; New = \"805\u00f77=\" },\n  @{ Old = \"688\u00f74=\"; New = \"854\u00f76=\" },\n  @{ Old = \"242\u00f73=\"; New = \"235\u00f79=\" },\n  @{ Old = \"700\u00f75=\"; New = \"511\u00f73=\" },\n  @{ Old = \"778\u00f78=\"; New = \"267\u00f75=\" },\n  @{ Old = \"702\u00f76=\"; New = \"525\u00f74=\" },\n  @{ Old = \"102\u00f78=\"; New = \"347\u00f79=\" },\n  @{ Old = \"289\u00f72=\"; New = \"245\u00f72=\" },\n  @{ Old = \"404\u00f75=\"; New = \"422\u00f76=\" },\n  @{ Old = \"947\u00f78=\"; New = \"314\u00f74=\" },\n  @{ Old = \"214\u00f78=\"; New = \"556\u00f76=\" },\n  @{ Old = \"491\u00f73=\"; New = \"324\u00f72=\" },\n  @{ Old = \"685\u00f75=\"; New = \"222\u00f79=\" },\n  @{ Old = \"960\u00f76=\"; New = \"431\u00f78=\" },\n  @{ Old = \"715\u00f76=\"; New = \"462\u00f77=\" },\n  @{ Old = \"868\u00f72=\"; New = \"401\u00f76=\" },\n  @{ Old = \"883\u00f74=\"; New = \"355\u00f77=\" },\n  @{ Old = \"438\u00f75=\"; New = \"839\u00f76=\" },\n  @{ Old = \"223\u00f76=\"; New = \"503\u00f74=\" },\n  @{ Old = \"713\u00f74=\"; New = \"659\u00f72=\" }\n)\n\nforeach ($pair in $replacements) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute(\n    $pair.Old,\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    $pair.New,\n    2\n  )\n}\n"}
